# Daily attendance processing - 2026-01-03 09:57:23
# Swap the order of the names listed in the "Recorded By" column (G) so that
# "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
